# Auto-upload VRF Excel file
# Adds a new worksheet named "fgh" at the end of the workbook, containing
# the standard VRF model header row (same header template used by every
# other sheet in this workbook).

$wb = $excel.ActiveWorkbook

# Worksheets.Add() inserts the new sheet right before the active sheet,
# so create it first, name it, then move it to be the very last tab.
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "fgh"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Move($null, $lastSheet)

# Re-fetch the sheet by name so subsequent writes land on the (now moved)
# worksheet object.
$ws = $wb.Worksheets.Item("fgh")

$ws.Range("A1").Value = "Outdoor Model"
$ws.Range("B1").Value = "Outdoor Quantity"
$ws.Range("C1").Value = "Outdoor Serial(s)"
$ws.Range("D1").Value = "Indoor Model"
$ws.Range("E1").Value = "Indoor Quantity"
$ws.Range("F1").Value = "Indoor Serial(s)"

# Match the bold / thin-bordered / center+top-aligned header formatting
# used by the other VRF sheets in this workbook.
$headerRange = $ws.Range("A1:F1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders.LineStyle = 1         # xlContinuous
$headerRange.Borders.Weight = 2            # xlThin
